$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5800840258598328
$ws.Range("B1").Value = 0.3779819905757904
$ws.Range("C1").Value = 0.2940379679203033
$ws.Range("D1").Value = 0.2856670916080475
$ws.Range("E1").Value = 0.3092442154884338
